$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 358) holds the "Förändrad" (last changed) date
# as a serial date number. Bump every value in that range from 45179 to
# 45180 (i.e. one day later), matching the source data refresh.
$ws.Range("C2:C358").Value = 45180
